$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.065.47"
$ws.Range("D3").Value = "1.726.48"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'219.00"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").Value = "'0.524"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +13.22%  "
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("D10").Value = "'0.0633"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "1.970.39"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("D13").Value = "1.733.76"
$ws.Range("E13").Value = "  +3.41%  "
$ws.Range("E15").Value = "  +5.30%  "
$ws.Range("D16").Value = "'67.55"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").Value = "28.035.69"
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").Value = "'243.26"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").Value = "0.0₃0755"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("E20").Value = "  -3.22%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "'4.63"
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("D23").Value = "'9.74"
$ws.Range("E23").Value = "  +4.29%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'148.89"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "'7.53"
$ws.Range("E26").Value = "  +4.29%  "
$ws.Range("D27").Value = "'16.76"
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("D32").Value = "'3.45"
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("D33").Value = "1.492.89"
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("D34").Value = "'3.27"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("D35").Value = "'1.66"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("D36").Value = "'0.954"
$ws.Range("E36").Value = "  +3.16%  "
$ws.Range("D37").Value = "'0.607"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").Value = "'70.79"
$ws.Range("E41").Value = "  +4.63%  "
$ws.Range("E42").Value = "  +4.13%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("D45").Value = "1.874.09"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").Value = "'0.797"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("E47").Value = "  +12.03%  "
$ws.Range("D48").Value = "'91.14"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  +4.12%  "
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").Value = "'8.19"
$ws.Range("E51").Value = "  +2.07%  "
